# =========================================================================
# Daily attendance processing - 2025-11-17 18:29:23
# Adds 13 new ANATOMY Session 3 attendance records to the Attendance sheet
# and refreshes the derived Summary-sheet statistics for the affected
# students (Total Attended/Missed, ANATOMY attendance, percentage, status).
# =========================================================================

$wb = $excel.ActiveWorkbook
$wsAttendance = $wb.Worksheets.Item("Attendance")
$wsSummary = $wb.Worksheets.Item("Summary")

# ---- 1) Append new raw attendance log rows (ANATOMY, Session 3, 17/11/2025) ----
# Row 563: student 220766
$wsAttendance.Cells.Item(563,1).Value() = '''220766'
$wsAttendance.Cells.Item(563,2).Value() = 'ميار بنت خالد بن محمد الشيخ'
$wsAttendance.Cells.Item(563,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(563,4).Value() = 'C1'
$wsAttendance.Cells.Item(563,5).Value() = '220766@med.asu.edu.eg'
$wsAttendance.Cells.Item(563,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(563,7).Value() = '''3'
$wsAttendance.Cells.Item(563,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(563,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(563,10).Value() = '10:16:02'
$wsAttendance.Cells.Item(563,11).Value() = 'C1'

# Row 564: student 212442
$wsAttendance.Cells.Item(564,1).Value() = '''212442'
$wsAttendance.Cells.Item(564,2).Value() = 'رميساء محى الدين الامين الطيب'
$wsAttendance.Cells.Item(564,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(564,4).Value() = 'C1'
$wsAttendance.Cells.Item(564,5).Value() = '212442@med.asu.edu.eg'
$wsAttendance.Cells.Item(564,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(564,7).Value() = '''3'
$wsAttendance.Cells.Item(564,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(564,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(564,10).Value() = '10:16:08'
$wsAttendance.Cells.Item(564,11).Value() = 'C1'

# Row 565: student 220428
$wsAttendance.Cells.Item(565,1).Value() = '''220428'
$wsAttendance.Cells.Item(565,2).Value() = 'بسمله محمد عبد الحميد محمد'
$wsAttendance.Cells.Item(565,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(565,4).Value() = 'C1'
$wsAttendance.Cells.Item(565,5).Value() = '220428@med.asu.edu.eg'
$wsAttendance.Cells.Item(565,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(565,7).Value() = '''3'
$wsAttendance.Cells.Item(565,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(565,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(565,10).Value() = '10:16:19'
$wsAttendance.Cells.Item(565,11).Value() = 'C1'

# Row 566: student 221682
$wsAttendance.Cells.Item(566,1).Value() = '''221682'
$wsAttendance.Cells.Item(566,2).Value() = 'سرين حاج صدوق'
$wsAttendance.Cells.Item(566,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(566,4).Value() = 'C1'
$wsAttendance.Cells.Item(566,5).Value() = '221682@med.asu.edu.eg'
$wsAttendance.Cells.Item(566,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(566,7).Value() = '''3'
$wsAttendance.Cells.Item(566,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(566,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(566,10).Value() = '10:16:27'
$wsAttendance.Cells.Item(566,11).Value() = 'C1'

# Row 567: student 212318
$wsAttendance.Cells.Item(567,1).Value() = '''212318'
$wsAttendance.Cells.Item(567,2).Value() = 'مازن مصدق يس عبد اللطيف'
$wsAttendance.Cells.Item(567,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(567,4).Value() = 'C1'
$wsAttendance.Cells.Item(567,5).Value() = '212318@med.asu.edu.eg'
$wsAttendance.Cells.Item(567,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(567,7).Value() = '''3'
$wsAttendance.Cells.Item(567,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(567,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(567,10).Value() = '10:16:35'
$wsAttendance.Cells.Item(567,11).Value() = 'C1'

# Row 568: student 220743
$wsAttendance.Cells.Item(568,1).Value() = '''220743'
$wsAttendance.Cells.Item(568,2).Value() = 'مصطفى كرم سلامه سليمان'
$wsAttendance.Cells.Item(568,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(568,4).Value() = 'C1'
$wsAttendance.Cells.Item(568,5).Value() = '220743@med.asu.edu.eg'
$wsAttendance.Cells.Item(568,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(568,7).Value() = '''3'
$wsAttendance.Cells.Item(568,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(568,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(568,10).Value() = '10:16:43'
$wsAttendance.Cells.Item(568,11).Value() = 'C1'

# Row 569: student 212322
$wsAttendance.Cells.Item(569,1).Value() = '''212322'
$wsAttendance.Cells.Item(569,2).Value() = 'مهند حافظ عابدين الفاضل'
$wsAttendance.Cells.Item(569,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(569,4).Value() = 'C1'
$wsAttendance.Cells.Item(569,5).Value() = '212322@med.asu.edu.eg'
$wsAttendance.Cells.Item(569,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(569,7).Value() = '''3'
$wsAttendance.Cells.Item(569,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(569,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(569,10).Value() = '10:16:49'
$wsAttendance.Cells.Item(569,11).Value() = 'C1'

# Row 570: student 220304
$wsAttendance.Cells.Item(570,1).Value() = '''220304'
$wsAttendance.Cells.Item(570,2).Value() = 'احمد الكامل محمد عبدون عثمان'
$wsAttendance.Cells.Item(570,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(570,4).Value() = 'C1'
$wsAttendance.Cells.Item(570,5).Value() = '220304@med.asu.edu.eg'
$wsAttendance.Cells.Item(570,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(570,7).Value() = '''3'
$wsAttendance.Cells.Item(570,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(570,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(570,10).Value() = '10:16:56'
$wsAttendance.Cells.Item(570,11).Value() = 'C1'

# Row 571: student 212308
$wsAttendance.Cells.Item(571,1).Value() = '''212308'
$wsAttendance.Cells.Item(571,2).Value() = 'سحر محمد يوسف محمد'
$wsAttendance.Cells.Item(571,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(571,4).Value() = 'C1'
$wsAttendance.Cells.Item(571,5).Value() = '212308@med.asu.edu.eg'
$wsAttendance.Cells.Item(571,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(571,7).Value() = '''3'
$wsAttendance.Cells.Item(571,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(571,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(571,10).Value() = '10:17:04'
$wsAttendance.Cells.Item(571,11).Value() = 'C1'

# Row 572: student 212125
$wsAttendance.Cells.Item(572,1).Value() = '''212125'
$wsAttendance.Cells.Item(572,2).Value() = 'راما محمد الحاج محمد'
$wsAttendance.Cells.Item(572,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(572,4).Value() = 'C1'
$wsAttendance.Cells.Item(572,5).Value() = '212125@med.asu.edu.eg'
$wsAttendance.Cells.Item(572,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(572,7).Value() = '''3'
$wsAttendance.Cells.Item(572,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(572,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(572,10).Value() = '10:17:18'
$wsAttendance.Cells.Item(572,11).Value() = 'C1'

# Row 573: student 220618
$wsAttendance.Cells.Item(573,1).Value() = '''220618'
$wsAttendance.Cells.Item(573,2).Value() = 'غيداء مرزوق الرفاعى'
$wsAttendance.Cells.Item(573,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(573,4).Value() = 'C1'
$wsAttendance.Cells.Item(573,5).Value() = '220618@med.asu.edu.eg'
$wsAttendance.Cells.Item(573,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(573,7).Value() = '''3'
$wsAttendance.Cells.Item(573,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(573,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(573,10).Value() = '10:17:30'
$wsAttendance.Cells.Item(573,11).Value() = 'C1'

# Row 574: student 212160
$wsAttendance.Cells.Item(574,1).Value() = '''212160'
$wsAttendance.Cells.Item(574,2).Value() = 'المى عماد تركمانى'
$wsAttendance.Cells.Item(574,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(574,4).Value() = 'C1'
$wsAttendance.Cells.Item(574,5).Value() = '212160@med.asu.edu.eg'
$wsAttendance.Cells.Item(574,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(574,7).Value() = '''3'
$wsAttendance.Cells.Item(574,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(574,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(574,10).Value() = '10:17:41'
$wsAttendance.Cells.Item(574,11).Value() = 'C1'

# Row 575: student 220314
$wsAttendance.Cells.Item(575,1).Value() = '''220314'
$wsAttendance.Cells.Item(575,2).Value() = 'احمد ربيع قطب عبد المطلب بهوت'
$wsAttendance.Cells.Item(575,3).Value() = 'Year 2'
$wsAttendance.Cells.Item(575,4).Value() = 'C1'
$wsAttendance.Cells.Item(575,5).Value() = '220314@med.asu.edu.eg'
$wsAttendance.Cells.Item(575,6).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(575,7).Value() = '''3'
$wsAttendance.Cells.Item(575,8).Value() = 'ANATOMY'
$wsAttendance.Cells.Item(575,9).Value() = '17/11/2025'
$wsAttendance.Cells.Item(575,10).Value() = '10:21:42'
$wsAttendance.Cells.Item(575,11).Value() = 'C1'

# ---- 2) Re-apply the AutoFilter over the new, larger data range ----
$wsAttendance.AutoFilterMode = $false
$wsAttendance.Range("A1:K575").AutoFilter()

# ---- 3) Point the sheet-local _FilterDatabase defined name at the new range ----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$575"
    }
}

# ---- 4) Refresh the per-student Summary rows impacted by the new attendance ----
# Template cells (outside the edited set) used to copy the correct Status fill/style:
#   row 2  -> 'High Risk' style    row 4  -> 'Fail' style    row 22 -> 'Moderate Risk' style
# Row 51
$wsSummary.Cells.Item(51,7).Value() = '6.9%'
$wsSummary.Cells.Item(51,8).Value() = 20
$wsSummary.Cells.Item(51,12).Value() = 2
$wsSummary.Cells.Item(51,13).Value() = 8
$wsSummary.Cells.Item(51,15).Value() = 1
$wsSummary.Cells.Item(51,18).Value() = 1

# Row 54
$wsSummary.Cells.Item(2,6).Copy()
$wsSummary.Cells.Item(54,6).PasteSpecial(-4122)
$wsSummary.Cells.Item(54,6).Value() = 'High Risk'
$wsSummary.Cells.Item(54,7).Value() = '10.3%'
$wsSummary.Cells.Item(54,8).Value() = 19
$wsSummary.Cells.Item(54,12).Value() = 3
$wsSummary.Cells.Item(54,13).Value() = 7
$wsSummary.Cells.Item(54,15).Value() = 2
$wsSummary.Cells.Item(54,18).Value() = 1

# Row 62
$wsSummary.Cells.Item(62,7).Value() = '6.9%'
$wsSummary.Cells.Item(62,8).Value() = 20
$wsSummary.Cells.Item(62,12).Value() = 2
$wsSummary.Cells.Item(62,13).Value() = 8
$wsSummary.Cells.Item(62,15).Value() = 1
$wsSummary.Cells.Item(62,18).Value() = 1

# Row 63
$wsSummary.Cells.Item(22,6).Copy()
$wsSummary.Cells.Item(63,6).PasteSpecial(-4122)
$wsSummary.Cells.Item(63,6).Value() = 'Moderate Risk'
$wsSummary.Cells.Item(63,7).Value() = '17.2%'
$wsSummary.Cells.Item(63,8).Value() = 17
$wsSummary.Cells.Item(63,12).Value() = 5
$wsSummary.Cells.Item(63,13).Value() = 5
$wsSummary.Cells.Item(63,15).Value() = 3
$wsSummary.Cells.Item(63,18).Value() = 1

# Row 64
$wsSummary.Cells.Item(64,7).Value() = '20.7%'
$wsSummary.Cells.Item(64,8).Value() = 16
$wsSummary.Cells.Item(64,12).Value() = 6
$wsSummary.Cells.Item(64,13).Value() = 4
$wsSummary.Cells.Item(64,15).Value() = 3
$wsSummary.Cells.Item(64,18).Value() = 1

# Row 71
$wsSummary.Cells.Item(71,7).Value() = '13.8%'
$wsSummary.Cells.Item(71,8).Value() = 18
$wsSummary.Cells.Item(71,12).Value() = 4
$wsSummary.Cells.Item(71,13).Value() = 6
$wsSummary.Cells.Item(71,15).Value() = 2
$wsSummary.Cells.Item(71,18).Value() = 1

# Row 75
$wsSummary.Cells.Item(75,7).Value() = '13.8%'
$wsSummary.Cells.Item(75,8).Value() = 18
$wsSummary.Cells.Item(75,12).Value() = 4
$wsSummary.Cells.Item(75,13).Value() = 6
$wsSummary.Cells.Item(75,15).Value() = 1
$wsSummary.Cells.Item(75,18).Value() = 1

# Row 76
$wsSummary.Cells.Item(76,7).Value() = '13.8%'
$wsSummary.Cells.Item(76,8).Value() = 18
$wsSummary.Cells.Item(76,12).Value() = 4
$wsSummary.Cells.Item(76,13).Value() = 6
$wsSummary.Cells.Item(76,15).Value() = 3
$wsSummary.Cells.Item(76,18).Value() = 1

# Row 78
$wsSummary.Cells.Item(78,7).Value() = '13.8%'
$wsSummary.Cells.Item(78,8).Value() = 18
$wsSummary.Cells.Item(78,12).Value() = 4
$wsSummary.Cells.Item(78,13).Value() = 6
$wsSummary.Cells.Item(78,15).Value() = 2
$wsSummary.Cells.Item(78,18).Value() = 1

# Row 81
$wsSummary.Cells.Item(2,6).Copy()
$wsSummary.Cells.Item(81,6).PasteSpecial(-4122)
$wsSummary.Cells.Item(81,6).Value() = 'High Risk'
$wsSummary.Cells.Item(81,7).Value() = '10.3%'
$wsSummary.Cells.Item(81,8).Value() = 19
$wsSummary.Cells.Item(81,12).Value() = 3
$wsSummary.Cells.Item(81,13).Value() = 7
$wsSummary.Cells.Item(81,15).Value() = 2
$wsSummary.Cells.Item(81,18).Value() = 1

# Row 83
$wsSummary.Cells.Item(83,7).Value() = '3.4%'
$wsSummary.Cells.Item(83,8).Value() = 21
$wsSummary.Cells.Item(83,12).Value() = 1
$wsSummary.Cells.Item(83,13).Value() = 9
$wsSummary.Cells.Item(83,15).Value() = 1
$wsSummary.Cells.Item(83,18).Value() = 1

# Row 84
$wsSummary.Cells.Item(84,7).Value() = '20.7%'
$wsSummary.Cells.Item(84,8).Value() = 16
$wsSummary.Cells.Item(84,12).Value() = 6
$wsSummary.Cells.Item(84,13).Value() = 4
$wsSummary.Cells.Item(84,15).Value() = 2
$wsSummary.Cells.Item(84,18).Value() = 1

# Row 176
$wsSummary.Cells.Item(176,7).Value() = '20.7%'
$wsSummary.Cells.Item(176,8).Value() = 16
$wsSummary.Cells.Item(176,12).Value() = 6
$wsSummary.Cells.Item(176,13).Value() = 4
$wsSummary.Cells.Item(176,15).Value() = 3
$wsSummary.Cells.Item(176,18).Value() = 1

$excel.CutCopyMode = $false
Write-Host "Attendance + Summary sheets updated."